# Apply the commit's change:
#  - "mise en commentaires des paramètres de recherche qu'on n'utilise plus"
#  - Bump the Metadata "Date" property value
#  - Insert a new "Jurisdiction" property row (with empty value) right before
#    the "Description" row on the Metadata sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Date property value (row 8, column B)
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# 2. Insert a new row above the current "Description" row (row 11) so the
#    new row becomes row 11 and everything below shifts down by one.
$ws.Rows.Item(11).EntireRow.Insert()

# Copy the formatting from the row that is now just below (the shifted
# "Description" row, now row 12) so the new row matches the sheet's
# standard data-row style/border instead of the insert's default style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new "Jurisdiction" property / empty value pair.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
